# Update column F (dSF) values per repull of data / mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 5
$ws.Range("F6").Value = 1
$ws.Range("F8").Value = -4
$ws.Range("F9").Value = 2
$ws.Range("F10").Value = -3
$ws.Range("F11").Value = 4
$ws.Range("F12").Value = -3
$ws.Range("F13").Value = 2
